$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row arrived for "Feria Lagunitas de Puerto Montt / Granada".
# It belongs right after the header block of existing rows, at row 28, pushing the
# previous rows 28-44 down to 29-45 (dimension grows from A1:T44 to A1:T45).
$ws.Rows("28:28").Insert()

$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44680
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100104
$ws.Range("H28").Value = "Frutos de pepita"
$ws.Range("I28").Value = 100104001
$ws.Range("J28").Value = "Granada"
$ws.Range("K28").Value = "Wonderfull"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 80
$ws.Range("N28").Value = 14000
$ws.Range("O28").Value = 15000
$ws.Range("P28").Value = 14500
$ws.Range("Q28").Value = "$/caja 14 kilos empedrada"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 1036
$ws.Range("T28").Value = 14
